$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "y07EvltGLV9gR9HFDizh"
$ws.Range("D2").Value = "iS351pPKoJXdWFbMqmlE"
$ws.Range("E2").Value = "DZleUqFJpTOq4bh4kH57"
$ws.Range("F2").Value = "p3fg89xOJePF6heOQVDM"
$ws.Range("C3").Value = "LEvgGJFv3wX8cyojhIdn"
$ws.Range("D3").Value = "hYs064OntCcSskDoe29d"
$ws.Range("E3").Value = "0vm6lbGxskI5nvn49mOt"
$ws.Range("F3").Value = "eLecmvMrHuDfB3dmjPLI"
$ws.Range("C4").Value = "omXXaJnlYO8TnlcaR1bV"
$ws.Range("D4").Value = "L6RfCNE3qyUWhDwW0Z7P"
$ws.Range("E4").Value = "Z29KNPT07DyWttPhrVEz"
$ws.Range("F4").Value = "QbBz3cnWSlq1p3p1sDYD"
$ws.Range("C5").Value = "Jr5dn8P7XSkPQr70r65C"
$ws.Range("D5").Value = "PpLn10FJWh3fcjWJPDtV"
$ws.Range("E5").Value = "rg1k6HLOcOudEfEU0xbK"
$ws.Range("F5").Value = "BKrgoltCi57jqNlbVJXH"
$ws.Range("D6").Value = "853gANX2S8C1u1fcYqFp"
$ws.Range("E6").Value = "tGMkvqnQSc8jlFBxlnla"
$ws.Range("C8").Value = "dt5o0ZXq0lHipSEKE45M"
$ws.Range("D8").Value = "U0sCYYTyNCvbKioZGdeX"
$ws.Range("E8").Value = "xbYu1RT5bo0qoGvPwzdE"
$ws.Range("F8").Value = "4bGIwsaTQLMyz7XkmhZD"
$ws.Range("C9").Value = "OQjVlR3f6VeEXqVBtPj0"
$ws.Range("D9").Value = "nONZGiPUYu271V64D46e"
$ws.Range("E9").Value = "naiVMpLtq4ddpF5GD19a"
$ws.Range("F9").Value = "HDgCQ5CaESXkIytbwIFp"
$ws.Range("C10").Value = "McB4FI9csbPqlUrGalii"
$ws.Range("D10").Value = "jlpoQFLE2TxAAOsUPKmj"
$ws.Range("E10").Value = "jGw2euU9E1AKj3Yg9jJe"
$ws.Range("F10").Value = "V0I2Hs5IYMCcrKAIGcki"
$ws.Range("C11").Value = "koAak2Yp38mtWKRStVNX"
$ws.Range("D11").Value = "z6R3kfcgwuTM8bM8ltxe"
$ws.Range("E11").Value = "9gHKdfkJ4KhGG7quRo4V"
$ws.Range("F11").Value = "JTBc3A0tCOZbcrc8PFqC"
$ws.Range("D12").Value = "dmQai36eavG8ynY5QGJx"
$ws.Range("E12").Value = "UQjHafbeiUVGygJT4qlr"
$ws.Range("C14").Value = "lueGC7wwZ3azXw4HsMvd"
$ws.Range("D14").Value = "4p4kvLKB88R2Pn7t1xdB"
$ws.Range("E14").Value = "pIy6zoLXqWGGl8pcisDp"
$ws.Range("F14").Value = "7WzS2e4OqovYF60zsgaI"
$ws.Range("C15").Value = "FKLFL7zgFoSW1FD4u7Tt"
$ws.Range("D15").Value = "QrXsuut2Qx0mbKxD3941"
$ws.Range("E15").Value = "eMBRHhJHDbxY1B1RMtU2"
$ws.Range("F15").Value = "UEG4aKqB2xMtmq0PJORj"
$ws.Range("C16").Value = "LqnaQ7E1UGcyjLdpfI1J"
$ws.Range("D16").Value = "xJKjiDhtN80alPYfcBQp"
$ws.Range("E16").Value = "FH9xFUbb19oYvvw5G7cR"
$ws.Range("F16").Value = "WM0s11nvAbxUPCjRJoHx"
$ws.Range("C17").Value = "3vEs6iDWUDco29GHx8Bq"
$ws.Range("D17").Value = "0RgG34FQQ9XYFpaZc4JR"
$ws.Range("E17").Value = "MxKv76wMA1Di5gblZXPN"
$ws.Range("F17").Value = "mYm8FBX8nH1QHtccFczk"
$ws.Range("D18").Value = "0ggjPNQFG217jBQvkqUu"
$ws.Range("E18").Value = "pzEqOSDr9Zz7WZJrVODU"

$ws.Range("F17").Select()
